$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "778÷6=" "664÷6="
Replace-Text "576÷2=" "204÷4="
Replace-Text "309÷5=" "181÷8="
Replace-Text "825÷8=" "243÷8="
Replace-Text "265÷5=" "860÷2="
Replace-Text "913÷7=" "399÷4="
Replace-Text "300÷5=" "882÷3="
Replace-Text "529÷6=" "519÷6="
Replace-Text "132÷6=" "411÷9="
Replace-Text "573÷5=" "226÷3="
Replace-Text "207÷7=" "616÷5="
Replace-Text "485÷2=" "920÷4="
Replace-Text "402÷2=" "657÷7="
Replace-Text "844÷7=" "766÷2="
Replace-Text "429÷5=" "131÷4="
Replace-Text "786÷6=" "456÷9="
Replace-Text "234÷4=" "549÷8="
Replace-Text "176÷4=" "817÷8="
Replace-Text "389÷8=" "227÷9="
Replace-Text "630÷8=" "967÷6="
Replace-Text "469÷9=" "556÷2="
Replace-Text "411÷3=" "706÷7="
Replace-Text "638÷7=" "912÷6="
Replace-Text "200÷7=" "846÷9="
Replace-Text "830÷2=" "332÷3="
